$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("LP1912")
$ws.Cells.Item(2, 1).Value = "Última actualización: 11:35:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 155"

$ws.Cells.Item(41, 1).Value = "07:52:32"
$ws.Cells.Item(41, 3).Value = "17_ROMERO"
$ws.Cells.Item(41, 4).Value = 8

$ws.Cells.Item(42, 1).Value = "06:38:54"
$ws.Cells.Item(42, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(42, 4).Value = 82

$ws.Cells.Item(51, 3).Value = "15_ABASTO"

$ws.Cells.Item(52, 3).Value = "11_ETCHEVERRY"

$ws.Cells.Item(64, 3).Value = "23_HERNANDEZ"

$ws.Cells.Item(65, 3).Value = "215B_EL PATO"

$ws.Cells.Item(75, 1).Value = "08:40:59"
$ws.Cells.Item(75, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(75, 4).Value = 37

$ws.Cells.Item(76, 1).Value = "08:30:14"
$ws.Cells.Item(76, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(76, 4).Value = 47

$ws.Cells.Item(77, 1).Value = "08:52:33"
$ws.Cells.Item(77, 3).Value = "14_ABASTO"
$ws.Cells.Item(77, 4).Value = 25

$ws.Cells.Item(85, 1).Value = "08:40:59"
$ws.Cells.Item(85, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(85, 4).Value = 62

$ws.Cells.Item(86, 1).Value = "08:30:14"
$ws.Cells.Item(86, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(86, 4).Value = 72

$ws.Cells.Item(102, 3).Value = "14_ABASTO"

$ws.Cells.Item(103, 3).Value = "15_ABASTO"

$ws.Cells.Item(112, 1).Value = "10:56:01"
$ws.Cells.Item(112, 3).Value = "81_EL PELIGRO"
$ws.Cells.Item(112, 4).Value = 5

$ws.Cells.Item(113, 1).Value = "09:23:52"
$ws.Cells.Item(113, 3).Value = "10_OLMOS"
$ws.Cells.Item(113, 4).Value = 98

$ws.Cells.Item(118, 3).Value = "14_ABASTO"

$ws.Cells.Item(119, 3).Value = "15X38_ABASTO"

$ws.Cells.Item(126, 1).Value = "11:35:40"
$ws.Cells.Item(126, 4).Value = 6

$ws.Cells.Item(128, 1).Value = "11:35:40"
$ws.Cells.Item(128, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(128, 4).Value = 10

$ws.Cells.Item(129, 1).Value = "10:56:01"
$ws.Cells.Item(129, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(129, 4).Value = 49

$ws.Cells.Item(131, 1).Value = "10:07:51"
$ws.Cells.Item(131, 3).Value = "225_GOMEZ"
$ws.Cells.Item(131, 4).Value = 105

$ws.Cells.Item(132, 1).Value = "11:13:01"
$ws.Cells.Item(132, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(132, 4).Value = 39

$ws.Cells.Item(133, 1).Value = "11:35:40"
$ws.Cells.Item(133, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(133, 4).Value = 18

$ws.Cells.Item(134, 1).Value = "11:35:40"
$ws.Cells.Item(134, 2).Value = "11:53"
$ws.Cells.Item(134, 3).Value = "225_GOMEZ"
$ws.Cells.Item(134, 4).Value = 18

$ws.Cells.Item(135, 1).Value = "11:35:40"
$ws.Cells.Item(135, 2).Value = "11:58"
$ws.Cells.Item(135, 3).Value = "17_ROMERO"
$ws.Cells.Item(135, 4).Value = 23

$ws.Cells.Item(136, 1).Value = "11:35:40"
$ws.Cells.Item(136, 2).Value = "12:05"
$ws.Cells.Item(136, 4).Value = 30

$ws.Cells.Item(137, 2).Value = "12:06"
$ws.Cells.Item(137, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(137, 4).Value = 53

$ws.Cells.Item(138, 1).Value = "11:35:40"
$ws.Cells.Item(138, 3).Value = "15_ABASTO"
$ws.Cells.Item(138, 4).Value = 35

$ws.Cells.Item(139, 1).Value = "11:35:40"
$ws.Cells.Item(139, 2).Value = "12:10"
$ws.Cells.Item(139, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(139, 4).Value = 35

$ws.Cells.Item(140, 1).Value = "11:35:40"
$ws.Cells.Item(140, 2).Value = "12:17"
$ws.Cells.Item(140, 3).Value = "10_OLMOS"
$ws.Cells.Item(140, 4).Value = 42

$ws.Cells.Item(141, 1).Value = "11:35:40"
$ws.Cells.Item(141, 2).Value = "12:21"
$ws.Cells.Item(141, 4).Value = 46

$ws.Cells.Item(142, 2).Value = "12:22"
$ws.Cells.Item(142, 3).Value = "215C_EL PATO"
$ws.Cells.Item(142, 4).Value = 69

$ws.Cells.Item(143, 1).Value = "11:35:40"
$ws.Cells.Item(143, 2).Value = "12:31"
$ws.Cells.Item(143, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(143, 4).Value = 56

$ws.Cells.Item(144, 2).Value = "12:31"
$ws.Cells.Item(144, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(144, 4).Value = 78

$ws.Cells.Item(145, 1).Value = "11:35:40"
$ws.Cells.Item(145, 2).Value = "12:32"
$ws.Cells.Item(145, 3).Value = "14_ABASTO"
$ws.Cells.Item(145, 4).Value = 57

$ws.Cells.Item(146, 2).Value = "12:33"
$ws.Cells.Item(146, 3).Value = "14_ABASTO"
$ws.Cells.Item(146, 4).Value = 80

$ws.Cells.Item(147, 1).Value = "11:13:01"
$ws.Cells.Item(147, 2).Value = "12:33"
$ws.Cells.Item(147, 3).Value = "23_HERNANDEZ"
$ws.Cells.Item(147, 4).Value = 80

$ws.Cells.Item(148, 1).Value = "10:56:01"
$ws.Cells.Item(148, 2).Value = "12:34"
$ws.Cells.Item(148, 4).Value = 98

$ws.Cells.Item(149, 1).Value = "11:35:40"
$ws.Cells.Item(149, 2).Value = "12:34"
$ws.Cells.Item(149, 3).Value = "15_ABASTO"
$ws.Cells.Item(149, 4).Value = 59

$ws.Cells.Item(150, 1).Value = "11:35:40"
$ws.Cells.Item(150, 2).Value = "12:36"
$ws.Cells.Item(150, 3).Value = "27_EL RETIRO"
$ws.Cells.Item(150, 4).Value = 61

$ws.Cells.Item(151, 1).Value = "11:35:40"
$ws.Cells.Item(151, 2).Value = "12:47"
$ws.Cells.Item(151, 3).Value = "14_ABASTO"
$ws.Cells.Item(151, 4).Value = 72
$ws.Cells.Item(151, 5).Value = "LP1912"

$ws.Cells.Item(152, 1).Value = "11:35:40"
$ws.Cells.Item(152, 2).Value = "12:48"
$ws.Cells.Item(152, 3).Value = "15X38_ABASTO"
$ws.Cells.Item(152, 4).Value = 73
$ws.Cells.Item(152, 5).Value = "LP1912"

$ws.Cells.Item(153, 1).Value = "11:35:40"
$ws.Cells.Item(153, 2).Value = "12:48"
$ws.Cells.Item(153, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(153, 4).Value = 73
$ws.Cells.Item(153, 5).Value = "LP1912"

$ws.Cells.Item(154, 1).Value = "11:35:40"
$ws.Cells.Item(154, 2).Value = "13:02"
$ws.Cells.Item(154, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(154, 4).Value = 87
$ws.Cells.Item(154, 5).Value = "LP1912"

$ws.Cells.Item(155, 1).Value = "11:35:40"
$ws.Cells.Item(155, 2).Value = "13:03"
$ws.Cells.Item(155, 3).Value = "215C_EL PATO"
$ws.Cells.Item(155, 4).Value = 88
$ws.Cells.Item(155, 5).Value = "LP1912"

$ws.Cells.Item(156, 1).Value = "11:13:01"
$ws.Cells.Item(156, 2).Value = "13:03"
$ws.Cells.Item(156, 3).Value = "11_ETCHEVERRY"
$ws.Cells.Item(156, 4).Value = 110
$ws.Cells.Item(156, 5).Value = "LP1912"

$ws.Cells.Item(157, 1).Value = "11:35:40"
$ws.Cells.Item(157, 2).Value = "13:13"
$ws.Cells.Item(157, 3).Value = "16_SANTA ANA"
$ws.Cells.Item(157, 4).Value = 98
$ws.Cells.Item(157, 5).Value = "LP1912"

$ws.Cells.Item(158, 1).Value = "11:35:40"
$ws.Cells.Item(158, 2).Value = "13:17"
$ws.Cells.Item(158, 3).Value = "10_OLMOS"
$ws.Cells.Item(158, 4).Value = 102
$ws.Cells.Item(158, 5).Value = "LP1912"

$ws.Cells.Item(159, 1).Value = "11:35:40"
$ws.Cells.Item(159, 2).Value = "13:25"
$ws.Cells.Item(159, 3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(159, 4).Value = 110
$ws.Cells.Item(159, 5).Value = "LP1912"

$ws.Cells.Item(160, 1).Value = "11:35:40"
$ws.Cells.Item(160, 2).Value = "13:33"
$ws.Cells.Item(160, 3).Value = "215A_EL PATO"
$ws.Cells.Item(160, 4).Value = 118
$ws.Cells.Item(160, 5).Value = "LP1912"


$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Cells.Item(2, 1).Value = "Última actualización: 11:35:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 23"

$ws.Cells.Item(23, 1).Value = "11:35:40"
$ws.Cells.Item(23, 4).Value = 6

$ws.Cells.Item(25, 1).Value = "11:35:40"
$ws.Cells.Item(25, 4).Value = 46

$ws.Cells.Item(27, 1).Value = "11:35:40"
$ws.Cells.Item(27, 2).Value = "13:03"
$ws.Cells.Item(27, 3).Value = "215C_EL PATO"
$ws.Cells.Item(27, 4).Value = 88
$ws.Cells.Item(27, 5).Value = "LP1912"

$ws.Cells.Item(28, 1).Value = "11:35:40"
$ws.Cells.Item(28, 2).Value = "13:33"
$ws.Cells.Item(28, 3).Value = "215A_EL PATO"
$ws.Cells.Item(28, 4).Value = 118
$ws.Cells.Item(28, 5).Value = "LP1912"


$ws = $wb.Worksheets.Item("6203-6173")
$ws.Cells.Item(2, 1).Value = "Última actualización: 11:35:40"
$ws.Cells.Item(3, 1).Value = "Total filas: 20"

$ws.Cells.Item(19, 1).Value = "10:07:51"
$ws.Cells.Item(19, 3).Value = "215A_LA PLATA"
$ws.Cells.Item(19, 4).Value = 23

$ws.Cells.Item(20, 1).Value = "08:52:33"
$ws.Cells.Item(20, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(20, 4).Value = 98

$ws.Cells.Item(24, 1).Value = "11:35:40"
$ws.Cells.Item(24, 2).Value = "13:11"
$ws.Cells.Item(24, 3).Value = "215C_LA PLATA"
$ws.Cells.Item(24, 4).Value = 96
$ws.Cells.Item(24, 5).Value = "L6203"

$ws.Cells.Item(25, 1).Value = "11:35:40"
$ws.Cells.Item(25, 2).Value = "13:20"
$ws.Cells.Item(25, 3).Value = "215B_LP-P MOR-1 Y 57"
$ws.Cells.Item(25, 4).Value = 105
$ws.Cells.Item(25, 5).Value = "L6173"

